# Auto-generated COM-interop script applying scheduled market-data refresh
# to the Mateus_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 158.52632
$ws.Range("I12").Value = 158.52632
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 158.52632
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = 11.47368

$ws.Range("H40").Value = 6101807.5
$ws.Range("I40").Value = 2759.9285
$ws.Range("J40").Value = 9264277
$ws.Range("K40").Value = 2759.9285
$ws.Range("L40").Value = 9264277
$ws.Range("M40").Value = -2584.9285
$ws.Range("N40").Value = -9264627

$ws.Range("H112").Value = 3479.375
$ws.Range("I112").Value = 2395.3333
$ws.Range("J112").Value = 4129.8
$ws.Range("K112").Value = 7185.999899999999
$ws.Range("L112").Value = 12389.4
$ws.Range("M112").Value = -6077.999899999999
$ws.Range("N112").Value = -14605.4

$ws.Range("H129").Value = 4630261
$ws.Range("I129").Value = 710.875
$ws.Range("J129").Value = 41666664
$ws.Range("K129").Value = 2132.625
$ws.Range("L129").Value = 124999992
$ws.Range("M129").Value = 2867.375
$ws.Range("N129").Value = -125009992

$ws.Range("H137").Value = 2750.2727
$ws.Range("I137").Value = 2730.077
$ws.Range("J137").Value = 2779.4443
$ws.Range("K137").Value = 8190.231000000001
$ws.Range("L137").Value = 8338.332900000001
$ws.Range("M137").Value = -5640.231000000001
$ws.Range("N137").Value = -13438.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3889.3
$ws.Range("I45").Value = 2991.9285
$ws.Range("J45").Value = 4674.5
$ws.Range("K45").Value = 2991.9285
$ws.Range("L45").Value = 4674.5
$ws.Range("M45").Value = -2614.9285
$ws.Range("N45").Value = -5428.5

$ws.Range("H122").Value = 2824.3
$ws.Range("I122").Value = 2301.5
$ws.Range("J122").Value = 3608.5
$ws.Range("K122").Value = 6904.5
$ws.Range("L122").Value = 10825.5
$ws.Range("M122").Value = -4454.5
$ws.Range("N122").Value = -15725.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 830.9
$ws.Range("I22").Value = 788.625
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 788.625
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -615.625
$ws.Range("N22").Value = -1346

$ws.Range("H57").Value = 60000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 60000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 60000
$ws.Range("N57").Value = -61440

$ws.Range("H136").Value = 60000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 60000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 60000
$ws.Range("N136").Value = -70200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 40005100
$ws.Range("I31").Value = 76926770
$ws.Range("J31").Value = 6624.9165
$ws.Range("K31").Value = 76926770
$ws.Range("L31").Value = 6624.9165
$ws.Range("M31").Value = -76926475
$ws.Range("N31").Value = -7214.9165

$ws.Range("H34").Value = 40005100
$ws.Range("I34").Value = 76926770
$ws.Range("J34").Value = 6624.9165
$ws.Range("K34").Value = 76926770
$ws.Range("L34").Value = 6624.9165
$ws.Range("M34").Value = -76926568
$ws.Range("N34").Value = -7028.9165

$ws.Range("H58").Value = 11234.2
$ws.Range("I58").Value = 4257.4
$ws.Range("J58").Value = 13559.8
$ws.Range("K58").Value = 4257.4
$ws.Range("L58").Value = 13559.8
$ws.Range("M58").Value = -4054.4
$ws.Range("N58").Value = -13965.8

$ws.Range("H107").Value = 3085.25
$ws.Range("I107").Value = 938.8
$ws.Range("J107").Value = 6662.6665
$ws.Range("K107").Value = 938.8
$ws.Range("L107").Value = 6662.6665
$ws.Range("M107").Value = 981.2
$ws.Range("N107").Value = -10502.6665

$ws.Range("H108").Value = 23250
$ws.Range("I108").Value = 34500
$ws.Range("J108").Value = 12000
$ws.Range("K108").Value = 34500
$ws.Range("L108").Value = 12000
$ws.Range("M108").Value = -30660
$ws.Range("N108").Value = -19680

$ws.Range("H123").Value = 61750
$ws.Range("I123").Value = 48500
$ws.Range("J123").Value = 75000
$ws.Range("K123").Value = 48500
$ws.Range("L123").Value = 75000
$ws.Range("M123").Value = -43600
$ws.Range("N123").Value = -84800

$ws.Range("H134").Value = 3281.6538
$ws.Range("I134").Value = 1872.619
$ws.Range("J134").Value = 9199.6
$ws.Range("K134").Value = 5617.857
$ws.Range("L134").Value = 27598.8
$ws.Range("M134").Value = -3082.857
$ws.Range("N134").Value = -32668.8

$ws.Range("H136").Value = 11234.2
$ws.Range("I136").Value = 4257.4
$ws.Range("J136").Value = 13559.8
$ws.Range("K136").Value = 12772.2
$ws.Range("L136").Value = 40679.39999999999
$ws.Range("M136").Value = -10222.2
$ws.Range("N136").Value = -45779.39999999999

$ws.Range("H141").Value = 221965.06
$ws.Range("I141").Value = 26743.5
$ws.Range("J141").Value = 249853.86
$ws.Range("K141").Value = 26743.5
$ws.Range("L141").Value = 249853.86
$ws.Range("M141").Value = -21563.5
$ws.Range("N141").Value = -260213.86

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H12").Value = 1363
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1363
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 4089
$ws.Range("N12").Value = -4435

$ws.Range("H68").Value = 45458320
$ws.Range("I68").Value = 100004216
$ws.Range("J68").Value = 3406.3333
$ws.Range("K68").Value = 300012648
$ws.Range("L68").Value = 10218.9999
$ws.Range("M68").Value = -300011837
$ws.Range("N68").Value = -11840.9999

$ws.Range("H71").Value = 45458320
$ws.Range("I71").Value = 100004216
$ws.Range("J71").Value = 3406.3333
$ws.Range("K71").Value = 900037944
$ws.Range("L71").Value = 30656.9997
$ws.Range("M71").Value = -900033888
$ws.Range("N71").Value = -38768.9997

$ws.Range("H81").Value = 17491.5
$ws.Range("I81").Value = 650
$ws.Range("J81").Value = 34333
$ws.Range("K81").Value = 1950
$ws.Range("L81").Value = 102999
$ws.Range("M81").Value = -827
$ws.Range("N81").Value = -105245

$ws.Range("H84").Value = 17491.5
$ws.Range("I84").Value = 650
$ws.Range("J84").Value = 34333
$ws.Range("K84").Value = 5850
$ws.Range("L84").Value = 308997
$ws.Range("M84").Value = -234
$ws.Range("N84").Value = -320229

$ws.Range("H113").Value = 1636.5
$ws.Range("I113").Value = 478.66666
$ws.Range("J113").Value = 1903.6923
$ws.Range("K113").Value = 1435.99998
$ws.Range("L113").Value = 5711.0769
$ws.Range("M113").Value = 734.0000199999999
$ws.Range("N113").Value = -10051.0769

$ws.Range("H122").Value = 4858.905
$ws.Range("I122").Value = 3491.2856
$ws.Range("J122").Value = 5542.7144
$ws.Range("K122").Value = 31421.5704
$ws.Range("L122").Value = 49884.4296
$ws.Range("M122").Value = -28971.5704
$ws.Range("N122").Value = -54784.4296

$ws.Range("H132").Value = 29413582
$ws.Range("I132").Value = 500000000
$ws.Range("J132").Value = 1931.25
$ws.Range("K132").Value = 4500000000
$ws.Range("L132").Value = 17381.25
$ws.Range("M132").Value = -4499997470
$ws.Range("N132").Value = -22441.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3828
$ws.Range("I14").Value = 3357.6
$ws.Range("J14").Value = 5004
$ws.Range("K14").Value = 3357.6
$ws.Range("L14").Value = 5004
$ws.Range("M14").Value = -3189.6
$ws.Range("N14").Value = -5340

$ws.Range("H17").Value = 5297.7
$ws.Range("I17").Value = 8
$ws.Range("J17").Value = 52905
$ws.Range("K17").Value = 8
$ws.Range("L17").Value = 52905
$ws.Range("M17").Value = 160
$ws.Range("N17").Value = -53241

$ws.Range("H97").Value = 854.96295
$ws.Range("I97").Value = 852.5789
$ws.Range("J97").Value = 860.625
$ws.Range("K97").Value = 852.5789
$ws.Range("L97").Value = 860.625
$ws.Range("M97").Value = -356.5789
$ws.Range("N97").Value = -1852.625

$ws.Range("H122").Value = 2879.6
$ws.Range("I122").Value = 3399
$ws.Range("J122").Value = 2749.75
$ws.Range("K122").Value = 10197
$ws.Range("L122").Value = 8249.25
$ws.Range("M122").Value = -7747
$ws.Range("N122").Value = -13149.25

$ws.Range("H126").Value = 4977.5835
$ws.Range("I126").Value = 4341.375
$ws.Range("J126").Value = 6250
$ws.Range("K126").Value = 13024.125
$ws.Range("L126").Value = 18750
$ws.Range("M126").Value = -10554.125
$ws.Range("N126").Value = -23690

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4422.077
$ws.Range("I7").Value = 4422.077
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4422.077
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4310.077

$ws.Range("H40").Value = 7044.8
$ws.Range("I40").Value = 6089.6
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 6089.6
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -5953.6
$ws.Range("N40").Value = -8272

$ws.Range("H46").Value = 5579.263
$ws.Range("I46").Value = 6684.5806
$ws.Range("J46").Value = 684.2857
$ws.Range("K46").Value = 6684.5806
$ws.Range("L46").Value = 684.2857
$ws.Range("M46").Value = -6496.5806
$ws.Range("N46").Value = -1060.2857

$ws.Range("H59").Value = 18250
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 18250
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 18250
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -19558

$ws.Range("H68").Value = 4947.4
$ws.Range("I68").Value = 3699.4
$ws.Range("J68").Value = 6195.4
$ws.Range("K68").Value = 3699.4
$ws.Range("L68").Value = 6195.4
$ws.Range("M68").Value = -2950.4
$ws.Range("N68").Value = -7693.4

$ws.Range("H71").Value = 4947.4
$ws.Range("I71").Value = 3699.4
$ws.Range("J71").Value = 6195.4
$ws.Range("K71").Value = 18497
$ws.Range("L71").Value = 30977
$ws.Range("M71").Value = -14753
$ws.Range("N71").Value = -38465

$ws.Range("H126").Value = 4422.077
$ws.Range("I126").Value = 4422.077
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13266.231
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -10796.231

$ws.Range("H132").Value = 25671.5
$ws.Range("I132").Value = 24184.846
$ws.Range("J132").Value = 44998
$ws.Range("K132").Value = 72554.538
$ws.Range("L132").Value = 134994
$ws.Range("M132").Value = -70024.538
$ws.Range("N132").Value = -140054

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 80000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 80000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180

$ws.Range("H111").Value = 81473.664
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 81473.664
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 81473.664
$ws.Range("N111").Value = -89653.664

$ws.Range("H126").Value = 7338.3
$ws.Range("I126").Value = 6547.25
$ws.Range("J126").Value = 10502.5
$ws.Range("K126").Value = 19641.75
$ws.Range("L126").Value = 31507.5
$ws.Range("M126").Value = -17171.75
$ws.Range("N126").Value = -36447.5

$ws.Range("H132").Value = 3524.8096
$ws.Range("I132").Value = 3106.8057
$ws.Range("J132").Value = 6032.8335
$ws.Range("K132").Value = 9320.417099999999
$ws.Range("L132").Value = 18098.5005
$ws.Range("M132").Value = -6790.417099999999
$ws.Range("N132").Value = -23158.5005
